$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final list of ClassName values (header + 16 unique class names),
# replacing the old 37-row list with this trimmed/edited 17-row list.
$values = @(
    "ClassName",
    "Alakli Allotment",
    "Alkali Allotment",
    "Butler Butte Allotment",
    "Elliot Creek Allotment",
    "Fish Lake Allotment",
    "Hershberger Allotment",
    "High Cascade Ranger District",
    "High Cascades Ranger District",
    "Local",
    "Rogue River-Siskiyou National Forests National Forest",
    "Rouge River National Forest",
    "Siskiyou Mountains Ranger District",
    "Siskiyou Mountains Ranger District Ranger District",
    "Tiller Ranger District",
    "Unknown Pasture",
    "Woodruff Allotment"
)

# Clear out everything below the old data so no stray rows remain.
$ws.Range("A1:A37").ClearContents()

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
